$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the four new worksheets (in order) at the end of the workbook:
#   User_Growth, Monster_Template, Dungeon_List, Payment_Grade
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGrowth = $wb.Worksheets.Add($null, $lastSheet)
$wsGrowth.Name = "User_Growth"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMonster = $wb.Worksheets.Add($null, $lastSheet)
$wsMonster.Name = "Monster_Template"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDungeon = $wb.Worksheets.Add($null, $lastSheet)
$wsDungeon.Name = "Dungeon_List"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPayment = $wb.Worksheets.Add($null, $lastSheet)
$wsPayment.Name = "Payment_Grade"

# ---------------------------------------------------------------------------
# User_Growth sheet data
# ---------------------------------------------------------------------------
$wsGrowth.Cells.Item(1,1).Value = "Level"
$wsGrowth.Cells.Item(1,2).Value = "Base_HP"
$wsGrowth.Cells.Item(1,3).Value = "Base_ATK"
$wsGrowth.Cells.Item(1,4).Value = "Base_DEF"

$growthRows = @(
    @(1,   10,    1,    1),
    @(5,   100,   10,   10),
    @(10,  1000,  100,  100),
    @(20,  2000,  200,  200),
    @(30,  4000,  400,  400),
    @(40,  6000,  600,  600),
    @(50,  8000,  800,  800),
    @(60,  12000, 1200, 1200),
    @(70,  16000, 1600, 1600),
    @(80,  20000, 2000, 2000),
    @(90,  25000, 2500, 2500),
    @(100, 30000, 3000, 3000)
)

$r = 2
foreach ($row in $growthRows) {
    $wsGrowth.Cells.Item($r,1).Value = $row[0]
    $wsGrowth.Cells.Item($r,2).Value = $row[1]
    $wsGrowth.Cells.Item($r,3).Value = $row[2]
    $wsGrowth.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Monster_Template sheet data
# ---------------------------------------------------------------------------
$wsMonster.Cells.Item(1,1).Value = "Monster_Type"
$wsMonster.Cells.Item(1,2).Value = "ATK_Ratio"
$wsMonster.Cells.Item(1,3).Value = "DEF_Ratio"
$wsMonster.Cells.Item(1,4).Value = "HP_Ratio"
$wsMonster.Cells.Item(1,5).Value = "Note"

$wsMonster.Cells.Item(2,1).Value = "Normal"
$wsMonster.Cells.Item(2,2).Value = 0.5
$wsMonster.Cells.Item(2,3).Value = 1
$wsMonster.Cells.Item(2,4).Value = 0.3
$wsMonster.Cells.Item(2,5).Value = "일반 몬스터"

$wsMonster.Cells.Item(3,1).Value = "Elite"
$wsMonster.Cells.Item(3,2).Value = 0.83
$wsMonster.Cells.Item(3,3).Value = 1
$wsMonster.Cells.Item(3,4).Value = 1
$wsMonster.Cells.Item(3,5).Value = "정예 몬스터"

$wsMonster.Cells.Item(4,1).Value = "Boss"
$wsMonster.Cells.Item(4,2).Value = 1
$wsMonster.Cells.Item(4,3).Value = 1
$wsMonster.Cells.Item(4,4).Value = 2
$wsMonster.Cells.Item(4,5).Value = "보스 몬스터"

$wsMonster.Columns.Item(1).ColumnWidth = 13.625
$wsMonster.Columns.Item(2).ColumnWidth = 10
$wsMonster.Columns.Item(3).ColumnWidth = 9.875
$wsMonster.Columns.Item(4).ColumnWidth = 9.125
$wsMonster.Columns.Item(5).ColumnWidth = 11.625

# ---------------------------------------------------------------------------
# Dungeon_List sheet data
# ---------------------------------------------------------------------------
$wsDungeon.Cells.Item(1,1).Value = "Dungeon_Name"
$wsDungeon.Cells.Item(1,2).Value = "Unlock_Level"
$wsDungeon.Cells.Item(1,3).Value = "Monster_Type"
$wsDungeon.Cells.Item(1,4).Value = "Monster_Count"
$wsDungeon.Cells.Item(1,5).Value = "Target_Survival_Ratio"

$wsDungeon.Cells.Item(2,1).Value = "던전 1-1"
$wsDungeon.Cells.Item(2,2).Value = 1
$wsDungeon.Cells.Item(2,3).Value = "Normal"
$wsDungeon.Cells.Item(2,4).Value = 5
$wsDungeon.Cells.Item(2,5).Value = 1.2

$wsDungeon.Cells.Item(3,1).Value = "던전 1-5 (보스)"
$wsDungeon.Cells.Item(3,2).Value = 5
$wsDungeon.Cells.Item(3,3).Value = "Boss"
$wsDungeon.Cells.Item(3,4).Value = 1
$wsDungeon.Cells.Item(3,5).Value = 0.5

$wsDungeon.Cells.Item(4,1).Value = "던전 2-1"
$wsDungeon.Cells.Item(4,2).Value = 10
$wsDungeon.Cells.Item(4,3).Value = "Normal"
$wsDungeon.Cells.Item(4,4).Value = 8
$wsDungeon.Cells.Item(4,5).Value = 1.1000000000000001

$wsDungeon.Columns.Item(1).ColumnWidth = 15.875
$wsDungeon.Columns.Item(2).ColumnWidth = 12.75
$wsDungeon.Columns.Item(3).ColumnWidth = 13.625
$wsDungeon.Columns.Item(4).ColumnWidth = 15
$wsDungeon.Columns.Item(5).ColumnWidth = 20.25

# ---------------------------------------------------------------------------
# Payment_Grade sheet data
# ---------------------------------------------------------------------------
$wsPayment.Cells.Item(1,1).Value = "Grade"
$wsPayment.Cells.Item(1,2).Value = "Stat_Multiplier"
$wsPayment.Cells.Item(1,3).Value = "Note"

$wsPayment.Cells.Item(2,1).Value = "Free (무과금)"
$wsPayment.Cells.Item(2,2).Value = 1
$wsPayment.Cells.Item(2,3).Value = "기준점"

$wsPayment.Cells.Item(3,1).Value = "Light (소과금)"
$wsPayment.Cells.Item(3,2).Value = 1.2
$wsPayment.Cells.Item(3,3).Value = "20% 더 강함"

$wsPayment.Cells.Item(4,1).Value = "Heavy (헤비과금)"
$wsPayment.Cells.Item(4,2).Value = 1.6
$wsPayment.Cells.Item(4,3).Value = "60% 더 강함"

$wsPayment.Columns.Item(1).ColumnWidth = 16.625
$wsPayment.Columns.Item(2).ColumnWidth = 14.125
$wsPayment.Columns.Item(3).ColumnWidth = 12.25

# ---------------------------------------------------------------------------
# Selections on each new sheet (cosmetic, matches original author's cursor)
# ---------------------------------------------------------------------------
[void]$wsGrowth.Range("E21").Select()
[void]$wsMonster.Range("F5").Select()
[void]$wsDungeon.Range("F1").Select()
[void]$wsPayment.Range("B7").Select()

# ---------------------------------------------------------------------------
# Make Payment_Grade the active/selected tab (last one activated in source)
# ---------------------------------------------------------------------------
$wsPayment.Activate()
